# Insert a new data row for 2026/02/04 07:00 (ranking 201) right after the
# existing 2026/02/04 02:00 row (row 767), shifting all subsequent rows down
# by one. This matches the daily auto-push pattern of appending a new
# timestamp entry for the current day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 768, pushing rows 768..809 down to 769..810.
$ws.Rows.Item(768).Insert()

# Populate the newly inserted row with the new data point. Column A holds
# date-like text (not a real Excel date), so prefix with an apostrophe to
# keep Excel from auto-converting "2026/02/04" into a date serial value,
# then restore the default "Normal" style so it matches its sibling cells.
$ws.Cells.Item(768, 1).Value = "'2026/02/04"
$ws.Cells.Item(768, 1).Style = "Normal"
$ws.Cells.Item(768, 2).Value = "水"
$ws.Cells.Item(768, 3).Value = 7
$ws.Cells.Item(768, 4).Value = 201
